{"js": "// Commit: \"Write only the available phase out fields\"\n// The \"Substance details\" table has a header row with six columns:\n//   Substance | Baseline substance | Replacement Substance |\n//   Phase out MT | Phase out CO2 | Phase out ODP\n// Only the first three of those fields are actually populated/available,\n// so the last three (\"Phase out MT\", \"Phase out CO2\", \"Phase out ODP\")\n// header/data columns are dropped from the table (and its grid).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const t of tables.items) {\n  t.load(\"values\");\n}\nawait context.sync();\n\n// Find the table whose header row matches the \"Substance details\" table\n// (identified by its distinctive header text) instead of assuming a\n// fixed table index.\nlet target = null;\nfor (const t of tables.items) {\n  const header = t.values && t.values[0] ? t.values[0] : [];\n  if (\n    header.length >= 6 &&\n    header[0] === \"Substance\" &&\n    header[1] === \"Baseline substance\" &&\n    header[2] === \"Replacement Substance\" &&\n    header[3] === \"Phase out MT\" &&\n    header[4] === \"Phase out CO2\" &&\n    header[5] === \"Phase out ODP\"\n  ) {\n    target = t;\n    break;\n  }\n}\n\nif (target) {\n  // Remove the last three columns (Phase out MT / Phase out CO2 / Phase out ODP),\n  // leaving Substance / Baseline substance / Replacement Substance.\n  target.deleteColumns(3, 3);\n  await context.sync();\n}\n", "ps1": "# Commit: \"Write only the available phase out fields\"\n#\n# The \"Substance details\" table has a header row with six columns:\n#   Substance | Baseline substance | Replacement Substance |\n#   Phase out MT | Phase out CO2 | Phase out ODP\n# Only the first three of those fields are actually populated/available,\n# so the last three (\"Phase out MT\", \"Phase out CO2\", \"Phase out ODP\")\n# header/data columns are dropped from the table (and its grid).\n\n$d = $word.ActiveDocument\n\n# Locate the target table by its header text instead of assuming a fixed\n# table index.\n$target = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    if ($t.Columns.Count -ge 6) {\n        $h1 = $t.Cell(1,1).Range.Text.TrimEnd([char]13, [char]7)\n        $h2 = $t.Cell(1,2).Range.Text.TrimEnd([char]13, [char]7)\n        $h3 = $t.Cell(1,3).Range.Text.TrimEnd([char]13, [char]7)\n        $h4 = $t.Cell(1,4).Range.Text.TrimEnd([char]13, [char]7)\n        $h5 = $t.Cell(1,5).Range.Text.TrimEnd([char]13, [char]7)\n        $h6 = $t.Cell(1,6).Range.Text.TrimEnd([char]13, [char]7)\n        if ($h1 -eq \"Substance\" -and $h2 -eq \"Baseline substance\" -and `\n            $h3 -eq \"Replacement Substance\" -and $h4 -eq \"Phase out MT\" -and `\n            $h5 -eq \"Phase out CO2\" -and $h6 -eq \"Phase out ODP\") {\n            $target = $t\n            break\n        }\n    }\n}\n\nif ($target -ne $null) {\n    # Remove the last three columns (Phase out MT / Phase out CO2 / Phase out ODP),\n    # leaving Substance / Baseline substance / Replacement Substance.\n    # Delete from the right so earlier column indices stay valid.\n    $target.Columns.Item(6).Delete()\n    $target.Columns.Item(5).Delete()\n    $target.Columns.Item(4).Delete()\n}\n"}
